$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add I0 and IF headers, matching style of existing header cells (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-22
$data = @(
    @(7,7),
    @(6,6),
    @(8,8),
    @(8,8),
    @(7,7),
    @(9,9),
    @(6,6),
    @(7,7),
    @(7,7),
    @(6,7),
    @(7,7),
    @(7,7),
    @(6,7),
    @(7,8),
    @(5,5),
    @(6,6),
    @(7,7),
    @(7,7),
    @(6,7),
    @(6,6),
    @(4,4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
